$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 239. This shifts the existing rows 239-271
# down to 240-272, matching the diff (which shows every record from the
# old row 239 onward now living one row further down, plus one brand new
# record occupying the former row 239's position).
$ws.Rows.Item(239).Insert()

# Populate the newly inserted row 239 with the new weekly record.
$ws.Range("A239").Value = 8
$ws.Range("B239").Value = "Terminal La Palmera de La Serena"
$ws.Range("C239").Value = "Coquimbo"
$ws.Range("D239").Value = 45127
$ws.Range("E239").Value = 4
$ws.Range("F239").Value = 100112001
$ws.Range("G239").Value = "Berenjena"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 460
$ws.Range("K239").Value = 8000
$ws.Range("L239").Value = 9000
$ws.Range("M239").Value = 8500
$ws.Range("N239").Value = "$/caja 50 unidades"
$ws.Range("O239").Value = "Región de Arica y Parinacota"
$ws.Range("P239").Value = 170
$ws.Range("Q239").Value = 50
$ws.Range("R239").Value = "Hortaliza"
